# Refresh the correlation/p-value table on sheet "corr" with a new run's
# numbers, and refresh the significance-star (E column) annotations that
# accompany each p-value. Numeric literals are written in plain decimal
# (not scientific-notation) form because the COM script parser here does
# not accept an `E` exponent suffix on numeric literals; the underlying
# double value, and therefore the saved XML, is identical either way.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.049000000000000002
$ws.Range("C2").Value = 0.45100000000000001
$ws.Range("E2").ClearContents()

$ws.Range("B3").Value = -0.055
$ws.Range("C3").Value = 0.39600000000000002
$ws.Range("E3").ClearContents()

$ws.Range("B4").Value = 0.042999999999999997
$ws.Range("C4").Value = 0.50900000000000001
$ws.Range("E4").ClearContents()

$ws.Range("B5").Value = -0.0050000000000000001
$ws.Range("C5").Value = 0.94499999999999995
$ws.Range("E5").ClearContents()

$ws.Range("B6").Value = 0.11799999999999999
$ws.Range("C6").Value = 0.067000000000000004
$ws.Range("E6").ClearContents()

$ws.Range("B7").Value = 0.035000000000000003
$ws.Range("C7").Value = 0.58899999999999997
$ws.Range("E7").ClearContents()

$ws.Range("B8").Value = -0.122
$ws.Range("C8").Value = 0.06
$ws.Range("E8").ClearContents()

$ws.Range("B9").Value = -0.13100000000000001
$ws.Range("C9").Value = 0.042999999999999997
$ws.Range("E9").Value = "*"

$ws.Range("B10").Value = -0.19400000000000001
$ws.Range("C10").Value = 0.0030000000000000001
$ws.Range("E10").Value = "**"

$ws.Range("B11").Value = -0.14199999999999999
$ws.Range("C11").Value = 0.028000000000000001
$ws.Range("E11").Value = "*"

$ws.Range("B12").Value = 0.021000000000000001
$ws.Range("C12").Value = 0.748
$ws.Range("E12").ClearContents()

$ws.Range("B13").Value = 0.017000000000000001
$ws.Range("C13").Value = 0.78900000000000003
$ws.Range("E13").ClearContents()

$ws.Range("B14").Value = -0.20799999999999999
$ws.Range("C14").Value = 0.001
$ws.Range("E14").Value = "***"

$ws.Range("B15").Value = -0.153
$ws.Range("C15").Value = 0.017999999999999999
$ws.Range("E15").Value = "*"

$ws.Range("B16").Value = 0.017000000000000001
$ws.Range("C16").Value = 0.78900000000000003
$ws.Range("E16").ClearContents()

$ws.Range("B17").Value = 0.105
$ws.Range("C17").Value = 0.104
$ws.Range("E17").ClearContents()

$ws.Range("B18").Value = 0.014
$ws.Range("C18").Value = 0.83
$ws.Range("E18").ClearContents()

$ws.Range("B19").Value = 0.106
$ws.Range("C19").Value = 0.10299999999999999
$ws.Range("E19").ClearContents()

$ws.Range("B20").Value = 0.0070000000000000001
$ws.Range("C20").Value = 0.90900000000000003
$ws.Range("E20").ClearContents()

$ws.Range("B21").Value = -0.050999999999999997
$ws.Range("C21").Value = 0.433
$ws.Range("E21").ClearContents()

$ws.Range("B22").Value = -0.035999999999999997
$ws.Range("C22").Value = 0.58499999999999996
$ws.Range("E22").ClearContents()

$ws.Range("B23").Value = 0.012999999999999999
$ws.Range("C23").Value = 0.84599999999999997
$ws.Range("E23").ClearContents()

$ws.Range("B24").Value = -0.17899999999999999
$ws.Range("C24").Value = 0.0060000000000000001
$ws.Range("E24").Value = "**"

$ws.Range("B25").Value = -0.069000000000000006
$ws.Range("C25").Value = 0.28799999999999998
$ws.Range("E25").ClearContents()

$ws.Range("B26").Value = -0.20200000000000001
$ws.Range("C26").Value = 0.002
$ws.Range("E26").Value = "**"

$ws.Range("B27").Value = -0.17799999999999999
$ws.Range("C27").Value = 0.0060000000000000001
$ws.Range("E27").Value = "**"

$ws.Range("B28").Value = -0.156
$ws.Range("C28").Value = 0.016
$ws.Range("E28").Value = "*"

$ws.Range("B29").Value = -0.124
$ws.Range("C29").Value = 0.055
$ws.Range("E29").ClearContents()

$ws.Range("B30").Value = -0.113
$ws.Range("C30").Value = 0.081000000000000003
$ws.Range("E30").ClearContents()

$ws.Range("B31").Value = -0.075999999999999998
$ws.Range("C31").Value = 0.24099999999999999
$ws.Range("E31").ClearContents()

$ws.Range("B32").Value = -0.16400000000000001
$ws.Range("C32").Value = 0.010999999999999999
$ws.Range("E32").Value = "*"

$ws.Range("B33").Value = -0.16800000000000001
$ws.Range("C33").Value = 0.0089999999999999993
$ws.Range("E33").Value = "**"

$ws.Range("B34").Value = 0.058999999999999997
$ws.Range("C34").Value = 0.36499999999999999
$ws.Range("E34").ClearContents()

$ws.Range("B35").Value = 0.035000000000000003
$ws.Range("C35").Value = 0.59099999999999997
$ws.Range("E35").ClearContents()

$ws.Range("B36").Value = -0.02
$ws.Range("C36").Value = 0.76
$ws.Range("E36").ClearContents()

$ws.Range("B37").Value = 0.067000000000000004
$ws.Range("C37").Value = 0.30399999999999999
$ws.Range("E37").ClearContents()

# Update the active cell selection to match the saved view state
$ws.Range("F16").Select()
